# "Added HungaryFC Test data"
#
# Creates a new "Hungary" worksheet (cloned from "Turkey", which is the
# existing market sheet with the same A1:D16 layout) at the end of the
# workbook, fills in the Hungary-specific market name / part number, and
# moves the "active sheet" / selection state from Turkey onto the new
# Hungary tab.

$wb = $excel.ActiveWorkbook

$turkey = $wb.Worksheets.Item("Turkey")

# Clone the Turkey sheet and drop the copy right after it - this carries
# over layout, styles, merged cells, column widths, etc. for free.
$turkey.Copy($null, $turkey)
$hungary = $wb.Worksheets.Item($wb.Worksheets.Count)
$hungary.Name = "Hungary"

# Market-specific values for the new sheet.
$hungary.Range("B2").Value = "Hungary Market"
$hungary.Range("B4").Value = "NGC-4308/T3593/T3618"

# Widen column B a bit so the new, longer part number is visible.
$hungary.Columns("B:B").ColumnWidth = 20.27

# Turkey is no longer the selected/active tab - select it as a whole
# (mirrors a "select all" state) before moving focus away from it.
$turkey.Cells.Select() | Out-Null

# Hungary becomes the active sheet with cell E3 selected.
$hungary.Range("E3").Select() | Out-Null
$hungary.Activate()
